$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FrameworkSQA")

# Update the data values in row 2 (formulas in C2, F2, I2, L2 recalc automatically)
$ws.Range("B2").Value = 2
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 4
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 6
$ws.Range("J2").Value = 7
$ws.Range("K2").Value = 8

# Update the selected cell/active cell on the sheet view
$ws.Activate()
$ws.Range("L2").Select()
